$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 66897
$ws.Range("B2").Value = "Emanuel Pacheco"
$ws.Range("C2").Value = "Juridico"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45084
$ws.Range("G2").Value = 4356.48

# Row 3
$ws.Range("A3").Value = 5562
$ws.Range("B3").Value = "Lorenzo Novais"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45091
$ws.Range("G3").Value = 5476.47

# Row 4
$ws.Range("A4").Value = 3951
$ws.Range("B4").Value = "Diego Sousa"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 7135.03

# Row 5
$ws.Range("A5").Value = 18271
$ws.Range("B5").Value = "Danilo Ferreira"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Consulta medica"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 2414.01

# Row 6
$ws.Range("A6").Value = 84840
$ws.Range("B6").Value = "Luiza Pires"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Consulta medica"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45090
$ws.Range("G6").Value = 6854.8

# Row 7
$ws.Range("A7").Value = 44822
$ws.Range("B7").Value = "Valentim Rocha"
$ws.Range("C7").Value = "TI"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 3455.34

# Row 8
$ws.Range("A8").Value = 939
$ws.Range("B8").Value = "Srta. Maria Vitória Vieira"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Doenca"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45100
$ws.Range("G8").Value = 7535.22

# Row 9
$ws.Range("A9").Value = 77228
$ws.Range("B9").Value = "Luísa Rocha"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Consulta medica"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45089
$ws.Range("G9").Value = 3974.81

# Row 10
$ws.Range("A10").Value = 26892
$ws.Range("B10").Value = "Agatha Camargo"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45103
$ws.Range("G10").Value = 4727.89

# Row 11
$ws.Range("A11").Value = 50511
$ws.Range("B11").Value = "Marcos Vinicius da Costa"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45090
$ws.Range("G11").Value = 8614.83
